$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sourced resistors: swap Vishay Dale SMD parts for Stackpole through-hole parts.
$ws.Range("A10").Value = "Stackpole Electronics RSMF2JT100R "
$ws.Range("B10").Value = "Through Hole 100 ohm 2 watt resistor"
$ws.Range("H10").Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/RSMF2JT100R/RSMF2JT100RCT-ND/2021858"

$ws.Range("A11").Value = "Stackpole Electronics  CF12JT10R0 "
$ws.Range("B11").Value = "Through Hole 10 ohm ½ watt Resistor"
$ws.Range("H11").Value = "https://www.digikey.com/product-detail/en/stackpole-electronics-inc/CF12JT10R0/CF12JT10R0CT-ND/1830446"

# New unit prices / quantities for the resistors (formulas in D/E/F recalc automatically).
$ws.Range("C10").Value = 0.28
$ws.Range("D10").Value = 9

$ws.Range("C11").Value = 0.1

# D7 had picked up a duplicate "General" style (style index 5, identical to style 0);
# normalize it back to the plain default style like the rest of the column.
$ws.Range("D7").NumberFormat = "General"

# Cursor was left on E17 when the workbook was saved.
$ws.Range("E17").Select() | Out-Null
